# Auto-generated Excel COM-interop script to apply cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price text that Excel would otherwise
# auto-convert to a Number/Date when assigned via .Value, so the Price
# column cells being touched are pre-formatted as Text ("@") to preserve
# them exactly as the literal strings from the source diff.
$dPriceCells = @("D2","D3","D5","D6","D7","D8","D9","D11","D12","D13","D14","D15","D16","D17","D19","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D33","D35","D36","D37","D38","D39","D40","D43","D44","D45","D48","D50","D51")
foreach ($addr in $dPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "72.801.88"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "3.951.58"
$ws.Range("E3").Value = "  -2.44%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "604.87"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("D6").Value = "172.90"
$ws.Range("E6").Value = "  +12.19%  "
$ws.Range("D7").Value = "0.685"
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.792"
$ws.Range("E9").Value = "  +4.40%  "
$ws.Range("E10").Value = "  +7.16%  "
$ws.Range("D11").Value = "56.26"
$ws.Range("E11").Value = "  +4.28%  "
$ws.Range("D12").Value = "0.0000329"
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").Value = "11.63"
$ws.Range("E13").Value = "  +5.65%  "
$ws.Range("D14").Value = "4.579.45"
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "21.54"
$ws.Range("E15").Value = "  +3.98%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.947.01"
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").Value = "14.12"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").Value = "72.734.65"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "445.00"
$ws.Range("D22").Value = "4.81"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "95.88"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").Value = "3.33"
$ws.Range("E24").Value = "  -5.46%  "
$ws.Range("D25").Value = "14.21"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").Value = "4.24"
$ws.Range("E26").Value = "  -3.39%  "
$ws.Range("D27").Value = "11.26"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "10.45"
$ws.Range("E28").Value = "  -3.45%  "
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").Value = "5.90"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").Value = "35.97"
$ws.Range("E30").Value = "  -2.67%  "
$ws.Range("D31").Value = "7.97"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  +2.30%  "
$ws.Range("D33").Value = "49.90"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("E34").Value = "  -3.96%  "
$ws.Range("D35").Value = "0.0₃0998"
$ws.Range("E35").Value = "  +13.91%  "
$ws.Range("D36").Value = "69.22"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").Value = "634.12"
$ws.Range("E37").Value = "  -8.15%  "
$ws.Range("D38").Value = "0.429"
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("D39").Value = "3.43"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "3.24"
$ws.Range("E43").Value = "  +44.68%  "
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").Value = "10.65"
$ws.Range("E44").Value = "  -5.74%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0479"
$ws.Range("E45").Value = "  -3.14%  "
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").Value = "3.41"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("E49").Value = "  +5.99%  "
$ws.Range("D50").Value = "2.78"
$ws.Range("E50").Value = "  -18.07%  "
$ws.Range("D51").Value = "2.839.92"
$ws.Range("E51").Value = "  +1.22%  "
